$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 137 (shifts existing rows 137..199 down to 138..200)
$ws.Rows.Item(137).Insert()

# Populate the newly inserted row 137 with the new data point
$ws.Range("A137").Value = 7
$ws.Range("B137").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C137").Value = "Ñuble"
$ws.Range("D137").Value = 44510
$ws.Range("E137").Value = 16
$ws.Range("F137").Value = 100114013
$ws.Range("G137").Value = "Zanahoria"
$ws.Range("H137").Value = "Sin especificar"
$ws.Range("I137").Value = "Primera"
$ws.Range("J137").Value = 120
$ws.Range("K137").Value = 7500
$ws.Range("L137").Value = 8000
$ws.Range("M137").Value = 7750
$ws.Range("N137").Value = "$/saco 20 kilos"
$ws.Range("O137").Value = "Provincia de Diguillín"
$ws.Range("P137").Value = 388
$ws.Range("Q137").Value = 20
$ws.Range("R137").Value = "Hortaliza"

# Make sure the D137 cell keeps the date-style numeric formatting used by the
# rest of the column (same style as its neighbours, e.g. D136/D138).
$ws.Range("D137").NumberFormat = $ws.Range("D138").NumberFormat()
